$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing header cells (casing fix: metadata4Ing -> metadata4ing)
$ws.Range("D1").Value = "metadata4ing_IRI"
$ws.Range("E1").Value = "metadata4ing_DESC"

# Add new column F with header and per-row definition values
$ws.Range("F1").Value = "metadata4ing_DEF"
# Match the bold/centered/bordered header style used by the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(""Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'"", 'en')]"
$ws.Range("F3").Value = "['To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]', 'To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])']"
$ws.Range("F4").Value = "[]"
$ws.Range("F5").Value = "[locstr('A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.', 'en')]"
$ws.Range("F6").Value = "[]"
$ws.Range("F7").Value = "[]"
$ws.Range("F8").Value = "[]"
$ws.Range("F9").Value = "[locstr('An activity is something that occurs over a period of time and acts upon or with entities; it may include consuming, processing, transforming, modifying, relocating, using, or generating entities.', 'en'), locstr('Eine Aktivität ist etwas, das über einen Zeitraum und mit oder an Entitäten erfolgt; dazu kann das Verbrauchen, Verarbeiten, Transformieren, Modifizieren, Relokalisieren, Verwenden, oder die Genese von Entitäten gehören.', 'de')]"
$ws.Range("F10").Value = "[]"
$ws.Range("F11").Value = "[locstr('A specific representation of a dataset. A dataset might be available in multiple serializations that may differ in various ways, including natural language, media-type or format, schematic organization, temporal and spatial resolution, level of detail or profiles (which might specify any or all of the above).', 'en')]"
$ws.Range("F12").Value = "[locstr(""Abstract description of a method for analysis, generation and transformation of data and material objects like, e.g., 'Fourier Transform Infrared Spectroscopy' or 'Molecular Dynamics Simulation'. Note that while values for relevant parameters should be provided, m4i:Method does not refer to the actual execution of the process that is described (use m4i:ProcessingStep for that), but rather a description of the underlying principle."", 'en')]"
$ws.Range("F13").Value = "[]"
$ws.Range("F14").Value = "[]"
$ws.Range("F15").Value = "[]"
$ws.Range("F16").Value = "[]"
$ws.Range("F17").Value = "[locstr('Variable, i.e., a conventional that is employed for something to which values can be assigned (this may include different notations for the same variable including, e.g., computational representations); this is to be understood very generally, e.g., a variable may be an array or set of elements or have another kind of internal structure', 'en')]"
$ws.Range("F18").Value = "[]"
$ws.Range("F19").Value = "[locstr('Quantity value, i.e., a value that has a numerical magnitude and a physical unit', 'en')]"
